$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move the active selection from I20 to A3 -------------------------
$ws.Range("A3").Select() | Out-Null

# --- 2. Shrink the narrow "operator" columns (A:U and W:AMJ) -------------
# The stored column width is ColumnWidth + 5/6 (Excel's char->stored-width
# offset for this font), quantized to 1/6 steps by the engine, so
# 3.8333333333333335 lands on the nearest achievable width to the target
# 4.66396761133603 (stored width 4.666666666666667).
$ws.Columns("1:21").ColumnWidth = 3.8333333333333335
$ws.Columns("23:1025").ColumnWidth = 3.8333333333333335

# --- 3. Populate the previously-implicit blank row 19 ---------------------
# Row 20 keeps its row number (no shift) - we just need A19:C19 to exist as
# blank cells carrying the default ("Normal", style index 0) formatting,
# same as the rest of the sheet's untouched blank cells. Setting a Value or
# a named Style directly mints a brand-new style record in this engine, so
# instead copy the *formats only* from an existing default-styled blank
# cell (U1) - that reuses style index 0 exactly, with no stray style rows.
$ws.Range("U1").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
